$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
"74+23=",
"26+17=",
"9+58=",
"21+41=",
"34+37=",
"48-38=",
"9+35=",
"57+4=",
"62-59=",
"72-25=",
"46+22=",
"37+45=",
"77-69=",
"46+52=",
"38-10=",
"92-74=",
"65-35=",
"18+61=",
"18+32=",
"93-33=",
"81+12=",
"80-14=",
"3+20=",
"32+55=",
"45+33=",
"0+71=",
"41-17=",
"0+0=",
"12-1=",
"80-31=",
"75+17=",
"67-50=",
"9+76=",
"4+21=",
"13+10=",
"62-55=",
"58-46=",
"86-13=",
"15+58=",
"81-73=",
"10+35=",
"46-27=",
"10+60=",
"74-71=",
"4+8=",
"80-43=",
"69-52=",
"88-87=",
"61+36=",
"79-21=",
"45+44=",
"81-58=",
"62+8=",
"13+18=",
"3+40=",
"34+2=",
"10+13=",
"5+40=",
"65+3=",
"53+44=",
"91-26=",
"67-59=",
"94+5=",
"70-61=",
"23-20=",
"80-39=",
"37-15=",
"41+45=",
"95+0=",
"7-5=",
"77+4=",
"45+47=",
"35+38=",
"70+23=",
"99-93=",
"7+36=",
"76+17=",
"38+56=",
"15+6=",
"44+53=",
"91-19=",
"41+51=",
"7+30=",
"80-76=",
"96-39=",
"66-52=",
"96-21=",
"49+11=",
"49-21=",
"6+2=",
"64+30=",
"92-10=",
"72-3=",
"61-25=",
"32-28=",
"48-31=",
"73-2=",
"23+22=",
"85-25=",
"17+17="
)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
Write-Host "Done: $idx cells updated"